$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'263.03"
$ws.Range("E2").Value = "'1.03%"
$ws.Range("D3").Value = "'26.65"
$ws.Range("E3").Value = "'-2.22%"
$ws.Range("D4").Value = "'4.694"
$ws.Range("E4").Value = "'0.39%"
$ws.Range("D5").Value = "'0.06104"
$ws.Range("E5").Value = "'-1.01%"
$ws.Range("D6").Value = "'6.703"
$ws.Range("E6").Value = "'0.51%"
$ws.Range("D7").Value = "'0.8495"
$ws.Range("E7").Value = "'-0.14%"
$ws.Range("D8").Value = "'0.9066"
$ws.Range("E8").Value = "'-0.78%"
$ws.Range("D9").Value = "'0.1409"
$ws.Range("E9").Value = "'0.00%"
$ws.Range("D10").Value = "'0.05025"
$ws.Range("E10").Value = "'2.71%"
$ws.Range("D11").Value = "'0.07093"
$ws.Range("E11").Value = "'0.02%"
$ws.Range("D12").Value = "'0.03122"
$ws.Range("E12").Value = "'0.64%"
$ws.Range("D13").Value = "'0.09047"
$ws.Range("E13").Value = "'-0.11%"
$ws.Range("D14").Value = "'0.001534"
$ws.Range("E14").Value = "'-0.58%"
$ws.Range("D15").Value = "'0.0006177"
$ws.Range("E15").Value = "'-0.17%"
$ws.Range("D16").Value = "'0.005994"
$ws.Range("E16").Value = "'-1.83%"
$ws.Range("D17").Value = "'3.449"
$ws.Range("E17").Value = "'-0.05%"
$ws.Range("D18").Value = "'3.164"
$ws.Range("E18").Value = "'0.43%"
$ws.Range("E19").Value = "'-0.60%"
$ws.Range("D21").Value = "'0.1280"
$ws.Range("E21").Value = "'-1.44%"
$ws.Range("D22").Value = "'4.076"
$ws.Range("E22").Value = "'-0.51%"
$ws.Range("D23").Value = "'0.04245"
$ws.Range("E23").Value = "'-0.01%"
$ws.Range("D24").Value = "'0.001180"
$ws.Range("E24").Value = "'-3.00%"
$ws.Range("D25").Value = "'0.004058"
$ws.Range("E25").Value = "'6.89%"
$ws.Range("E26").Value = "'-0.02%"
$ws.Range("E27").Value = "'23.05%"
$ws.Range("D40").Value = "'0.03948"
$ws.Range("E40").Value = "'1.97%"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("D42").Value = "'0.004184"
$ws.Range("E42").Value = "'2.53%"
$ws.Range("D43").Value = "'0.002109"
$ws.Range("E43").Value = "'-3.95%"
$ws.Range("D44").Value = "'0.01158"
$ws.Range("E44").Value = "'-29.06%"
$ws.Range("D45").Value = "'0.00005095"
$ws.Range("E45").Value = "'-1.31%"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("D48").Value = "'0.2581"
$ws.Range("E48").Value = "'56.06%"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E50").Value = "'0.00%"
